$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each year block (rows laid out as A/B/C/D sub-periods), the B-row and
# C-row have had their contents (columns A:E) swapped.
for ($year = 2000; $year -le 2015; $year++) {
    $base = 2 + ($year - 2000) * 4
    $rowB = $base + 1
    $rowC = $base + 2

    $aB = $ws.Cells.Item($rowB, 1).Value2
    $bB = $ws.Cells.Item($rowB, 2).Value2
    $cB = $ws.Cells.Item($rowB, 3).Value2
    $dB = $ws.Cells.Item($rowB, 4).Value2
    $eB = $ws.Cells.Item($rowB, 5).Value2

    $aC = $ws.Cells.Item($rowC, 1).Value2
    $bC = $ws.Cells.Item($rowC, 2).Value2
    $cC = $ws.Cells.Item($rowC, 3).Value2
    $dC = $ws.Cells.Item($rowC, 4).Value2
    $eC = $ws.Cells.Item($rowC, 5).Value2

    $ws.Cells.Item($rowB, 1).Value = $aC
    $ws.Cells.Item($rowB, 2).Value = $bC
    $ws.Cells.Item($rowB, 3).Value = $cC
    $ws.Cells.Item($rowB, 4).Value = $dC
    $ws.Cells.Item($rowB, 5).Value = $eC

    $ws.Cells.Item($rowC, 1).Value = $aB
    $ws.Cells.Item($rowC, 2).Value = $bB
    $ws.Cells.Item($rowC, 3).Value = $cB
    $ws.Cells.Item($rowC, 4).Value = $dB
    $ws.Cells.Item($rowC, 5).Value = $eB
}

# Columns F (产销率) and G (销售量) — header + all data — are removed
# entirely, shrinking the used range from A1:G65 to A1:E65.
$ws.Range("F1:G65").Delete()
